$d = $word.ActiveDocument

function Replace-Text($oldText, $newText) {
    $d.Content.Find.Execute($oldText, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $newText, 2)
}

# Update the date in the title paragraph
Replace-Text "2025-08-05 Tuesday" "2025-08-06 Wednesday"

$tbl = $d.Tables.Item(1)

function Set-Cell($row, $col, $newText) {
    $cell = $tbl.Cell($row, $col)
    $rng = $cell.Range
    $rng.End = $rng.End - 1
    $rng.Text = $newText
}

# Row 1
Set-Cell 1 1 "335×8=2680"
Set-Cell 1 2 "810×2=1620"
Set-Cell 1 3 "664×2=1328"
Set-Cell 1 4 "596×5=2980"
Set-Cell 1 5 "369×2=738"

# Row 5
Set-Cell 5 1 "824×3=2472"
Set-Cell 5 2 "296×8=2368"
Set-Cell 5 3 "153×9=1377"
Set-Cell 5 4 "109×7=763"
Set-Cell 5 5 "435×8=3480"

# Row 10
Set-Cell 10 1 "946×3=2838"
Set-Cell 10 2 "899×7=6293"
Set-Cell 10 3 "612×2=1224"
Set-Cell 10 4 "322×8=2576"
Set-Cell 10 5 "810×8=6480"

# Row 15
Set-Cell 15 1 "965×5=4825"
Set-Cell 15 2 "250×7=1750"
Set-Cell 15 3 "627×9=5643"
Set-Cell 15 4 "232×9=2088"
Set-Cell 15 5 "225×8=1800"

# Row 20
Set-Cell 20 1 "279×4=1116"
Set-Cell 20 2 "704×6=4224"
Set-Cell 20 3 "267×7=1869"
Set-Cell 20 4 "866×7=6062"
Set-Cell 20 5 "742×7=5194"
